$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the existing row 476, pushing the rest of the
# table (old rows 476-558) down to 479-561.
$ws.Range("A476:A478").EntireRow.Insert()

# Row 476: new weekly entry - Murcott / Especial
$ws.Cells.Item(476,1).Value  = 4
$ws.Cells.Item(476,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(476,3).Value  = "Los Lagos"
$ws.Cells.Item(476,4).Value  = 45258
$ws.Cells.Item(476,5).Value  = 10
$ws.Cells.Item(476,6).Value  = "Fruta"
$ws.Cells.Item(476,7).Value  = 100102
$ws.Cells.Item(476,8).Value  = "Cítricos"
$ws.Cells.Item(476,9).Value  = 100102004
$ws.Cells.Item(476,10).Value = "Mandarina"
$ws.Cells.Item(476,11).Value = "Murcott"
$ws.Cells.Item(476,12).Value = "Especial"
$ws.Cells.Item(476,13).Value = 300
$ws.Cells.Item(476,14).Value = 10000
$ws.Cells.Item(476,15).Value = 10000
$ws.Cells.Item(476,16).Value = 10000
$ws.Cells.Item(476,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(476,18).Value = "Región de O'Higgins"
$ws.Cells.Item(476,19).Value = 1000
$ws.Cells.Item(476,20).Value = 10

# Row 477: new weekly entry - Murcott / Primera
$ws.Cells.Item(477,1).Value  = 4
$ws.Cells.Item(477,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(477,3).Value  = "Los Lagos"
$ws.Cells.Item(477,4).Value  = 45258
$ws.Cells.Item(477,5).Value  = 10
$ws.Cells.Item(477,6).Value  = "Fruta"
$ws.Cells.Item(477,7).Value  = 100102
$ws.Cells.Item(477,8).Value  = "Cítricos"
$ws.Cells.Item(477,9).Value  = 100102004
$ws.Cells.Item(477,10).Value = "Mandarina"
$ws.Cells.Item(477,11).Value = "Murcott"
$ws.Cells.Item(477,12).Value = "Primera"
$ws.Cells.Item(477,13).Value = 300
$ws.Cells.Item(477,14).Value = 8500
$ws.Cells.Item(477,15).Value = 8500
$ws.Cells.Item(477,16).Value = 8500
$ws.Cells.Item(477,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(477,18).Value = "Región de O'Higgins"
$ws.Cells.Item(477,19).Value = 850
$ws.Cells.Item(477,20).Value = 10

# Row 478: new weekly entry - Murcott / Segunda
$ws.Cells.Item(478,1).Value  = 4
$ws.Cells.Item(478,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(478,3).Value  = "Los Lagos"
$ws.Cells.Item(478,4).Value  = 45258
$ws.Cells.Item(478,5).Value  = 10
$ws.Cells.Item(478,6).Value  = "Fruta"
$ws.Cells.Item(478,7).Value  = 100102
$ws.Cells.Item(478,8).Value  = "Cítricos"
$ws.Cells.Item(478,9).Value  = 100102004
$ws.Cells.Item(478,10).Value = "Mandarina"
$ws.Cells.Item(478,11).Value = "Murcott"
$ws.Cells.Item(478,12).Value = "Segunda"
$ws.Cells.Item(478,13).Value = 300
$ws.Cells.Item(478,14).Value = 7000
$ws.Cells.Item(478,15).Value = 7000
$ws.Cells.Item(478,16).Value = 7000
$ws.Cells.Item(478,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(478,18).Value = "Región de O'Higgins"
$ws.Cells.Item(478,19).Value = 700
$ws.Cells.Item(478,20).Value = 10
